# Update "want to go" (F column) counts across the four sheets to match
# the newly generated data snapshot.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1240
$ws1.Range("F3").Value = 73
$ws1.Range("F5").Value = 3500
$ws1.Range("F6").Value = 1756
$ws1.Range("F7").Value = 6306
$ws1.Range("F8").Value = 136
$ws1.Range("F9").Value = 1893
$ws1.Range("F10").Value = 504
$ws1.Range("F11").Value = 9
$ws1.Range("F12").Value = 27
$ws1.Range("F15").Value = 47
$ws1.Range("F16").Value = 7504
$ws1.Range("F28").Value = 1704
$ws1.Range("F29").Value = 792
$ws1.Range("F30").Value = 358

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 359

# Sheet "本地生活" (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F4").Value = 673
$ws3.Range("F5").Value = 260

# Sheet "全部类型" (All Types) - aggregated view
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 673
$ws4.Range("F5").Value = 1240
$ws4.Range("F9").Value = 359
$ws4.Range("F10").Value = 3500
$ws4.Range("F11").Value = 260
$ws4.Range("F12").Value = 1756
$ws4.Range("F13").Value = 6306
$ws4.Range("F14").Value = 136
$ws4.Range("F15").Value = 1893
$ws4.Range("F17").Value = 504
$ws4.Range("F18").Value = 9
$ws4.Range("F19").Value = 27
$ws4.Range("F22").Value = 47
$ws4.Range("F23").Value = 7504
$ws4.Range("F34").Value = 1704
$ws4.Range("F35").Value = 792
$ws4.Range("F37").Value = 358

$wb.Save()
